$d = $word.ActiveDocument

# ===========================================================================
# Helper: stamp a brand-new run at $insertPos by pasting a copy of
# $donorRange (length $donorLen, already Copy()-ed format) and then
# overwriting its text with $newText. This reproduces the donor's *exact*
# run-properties (including the *absence* of a <w:color> element when the
# donor has none) because Word's Paste keeps the source formatting, whereas
# assigning Font.Color always re-serialises an explicit value (even "auto").
# Returns the end offset of the freshly written text.
# ===========================================================================
function Paste-Run($insertPos, $donorRange, $donorLen, $newText) {
    $donorRange.Copy()
    $ip = $d.Range($insertPos, $insertPos)
    $ip.Paste()
    $pasted = $d.Range($insertPos, $insertPos + $donorLen)
    $pasted.Text = $newText
    return $insertPos + $newText.Length
}

# ===========================================================================
# 1) Simple text-only fixes: drop the period after "lb" in two places.
# ===========================================================================
$d.Content.Find.Execute(" twenty or xxv lb. of ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " twenty or xxv lb of ", 2)

$d.Content.Find.Execute(" a quarter lb. of powder, half an ounce of glass is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " a quarter lb of powder, half an ounce of glass is ", 2)

# ===========================================================================
# 2) ". The metal of a big bell is made with 3 quintals of <del>one</del> <m>rosette"
#    -> ". The metal of a big bell is made with <del>one</del> <add>3</add> quintals of <m>rosette"
#    "one" loses its black (000000) run color, and the new "<add>"/"add>"
#    marker text is colorless too, while "3" and " quintals of " keep black.
# ===========================================================================

# 2a) Trim the "3 quintals of " text that used to precede "<del>one</del>".
$d.Content.Find.Execute(". The metal of a big bell is made with 3 quintals of ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". The metal of a big bell is made with ", 2)

# 2b) Donor carrying the colorless rPr (<w:rtl w:val="0"/> only): the word
#     "bigger" inside "<del>bigger</del>" earlier in the same paragraph.
$donorFind = $d.Content
$donorFind.Find.Execute("<del>bigger</del>")
$colorlessDonor = $d.Range($donorFind.Start + 5, $donorFind.Start + 5 + 6)

# 2c) Donor carrying the black (000000) rPr: the space in " clearer sound",
#     immediately before our edit area, left completely untouched.
$donor2Find = $d.Content
$donor2Find.Find.Execute(" clearer sound")
$blackDonor = $d.Range($donor2Find.Start, $donor2Find.Start + 1)

# 2d) Replace "one" (currently black) with a colorless copy of the same text.
$target = $d.Content
$target.Find.Execute("<del>one</del>")
$oneStart = $target.Start + 5
$oneEnd = $target.End - 6
$d.Range($oneStart, $oneEnd).Text = ""
$pos = Paste-Run $oneStart $colorlessDonor 6 "one"

# 2e) The single space that used to follow "</del>" becomes " <add>" and
#     loses its black color too. "</del>" is 6 chars.
$pos = $pos + 6
$d.Range($pos, $pos + 1).Text = ""
$pos = Paste-Run $pos $colorlessDonor 6 " <add>"

# 2f) New run "3</" keeps the black color.
$pos = Paste-Run $pos $blackDonor 1 "3</"

# 2g) New run "add>" is colorless.
$pos = Paste-Run $pos $colorlessDonor 6 "add>"

# 2h) New run " quintals of " keeps the black color.
$pos = Paste-Run $pos $blackDonor 1 " quintals of "
